$d = $word.ActiveDocument

$replacements = @(
    @{old = "PRINCIPAL CONSULTANT - Siege Analytics, Washington, DC | January 2014 – Present"; new = "PARTNER - Siege Analytics, Washington, DC | January 2014 – Present"},
    @{old = "PRINCIPAL CONSULTANT - Clarity and Rigour, Washington, DC | 2012 – 2014"; new = "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014"},
    @{old = "SENIOR CONSULTANT - Helm, Washington, DC | 2010 – 2012"; new = "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012"},
    @{old = "CONSULTANT - GSD&M, Austin, TX | 2008 – 2010"; new = "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010"},
    @{old = "SENIOR CONSULTANT - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008"; new = "RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008"},
    @{old = "CONSULTANT - Salsa Labs, Inc., Washington, DC | 2004 – 2006"; new = "SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006"},
    @{old = "CONSULTANT - The Praxis Project, Oakland, CA | 2002 – 2004"; new = "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004"},
    @{old = "ANALYST - Lake Research Partners, Washington, DC | 2001 – 2002"; new = "PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002"},
    @{old = "FIELD COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001"; new = "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
